# Update bank statement statement_135.xlsx content to new values per diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Name
$ws.Range("C2").Value = "Hartmut"

# Row 3 - Card number and surname
# (B3 holds a long digit string that must stay textual rather than be
#  auto-converted to a number; format as Text, assign it, then restore
#  the original cell formatting by pasting formats from an untouched
#  neighbor that already carries the same style.)
$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value = "2570314725427075"
$ws.Range("D3").Copy() | Out-Null
$ws.Range("B3").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("C3").Value = "Mohaupt"

# Row 5 - Opening balance date
$ws.Range("D5").Value = "KONTOSTAND AM 26.10.2023"

# Row 6 - transaction 1
$ws.Range("B6").Value = "29.10."
$ws.Range("C6").Value = "30.10."
$ws.Range("D6").Value = "BURGER KING Vilsbiburg"
$ws.Range("E6").Value = "41,83-"

# Row 7 - transaction 2
$ws.Range("B7").Value = "31.10."
$ws.Range("C7").Value = "01.11."
$ws.Range("D7").Value = "KARTENZAHLUNG JET TANKSTELLE"
$ws.Range("E7").Value = "84,47-"

# Row 8 - transaction 3
$ws.Range("B8").Value = "02.11."
$ws.Range("C8").Value = "03.11."
$ws.Range("D8").Value = "RECHNUNG VODAFONE GMBH 47074497"
$ws.Range("E8").Value = "39,16-"

# Row 9 - clear (formerly transaction 4); amount cell alignment changes to center
$ws.Range("B9").Value = ""
$ws.Range("C9").Value = ""
$ws.Range("D9").Value = ""
$ws.Range("E9").Value = ""
$ws.Range("E9").HorizontalAlignment = -4108  # xlCenter
$ws.Range("E9").VerticalAlignment = -4108    # xlCenter
$ws.Range("E9").WrapText = $true

# Row 10 - clear (formerly transaction 5); amount cell alignment changes to right/center/wrap
$ws.Range("B10").Value = ""
$ws.Range("C10").Value = ""
$ws.Range("D10").Value = ""
$ws.Range("E10").Value = ""
$ws.Range("E10").HorizontalAlignment = -4152  # xlRight
$ws.Range("E10").VerticalAlignment = -4108    # xlCenter
$ws.Range("E10").WrapText = $true

# Row 11 - clear (formerly transaction 6); amount cell alignment changes to right/center/wrap
$ws.Range("B11").Value = ""
$ws.Range("C11").Value = ""
$ws.Range("D11").Value = ""
$ws.Range("E11").Value = ""
$ws.Range("E11").HorizontalAlignment = -4152  # xlRight
$ws.Range("E11").VerticalAlignment = -4108    # xlCenter
$ws.Range("E11").WrapText = $true

# Row 12 - Closing balance date and amount
$ws.Range("D12").Value = "KONTOSTAND AM 05.11.2023"
$ws.Range("E12").Value = "165,46-"

# Row 13 - Next billing date
$ws.Range("C13").Value = "IHR NAECHSTER ABRECHNUNGSTERMIN 13.11.2023"
